$wb = $excel.ActiveWorkbook

# --- 1. Create the new "setting" sheet and place it before the existing "남자" sheet ---
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("남자"))
$newSheet.Name = "setting"

$settingSheet = $wb.Worksheets.Item("setting")
$settingSheet.Range("A1").Value = "sheet name"
$settingSheet.Range("B1").Value = "folder name"
$settingSheet.Range("A2").Value = 1
$settingSheet.Range("B2").Value = "analysis_man"
$settingSheet.Range("A3").Value = 2
$settingSheet.Range("B3").Value = "analysis_woman"

$settingSheet.Columns.Item(2).ColumnWidth = 11.714285714285714

# --- 2. Rename the existing sheets ("남자" -> "1", "여자" -> "2") ---
$wb.Worksheets.Item("남자").Name = "1"
$wb.Worksheets.Item("여자").Name = "2"

# --- 3. Append ".jpg" to every celebrity name in column B of the women's sheet ("2") ---
$womenSheet = $wb.Worksheets.Item("2")
$lastRow = $womenSheet.Cells($womenSheet.Rows.Count, 2).End(-4162).Row
for ($r = 3; $r -le $lastRow; $r++) {
    $cell = $womenSheet.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -ne $null -and -not ($val.ToString().EndsWith(".jpg"))) {
        $cell.Value = $val.ToString() + ".jpg"
    }
}

# --- 4. Cosmetic tweaks to match the authored workbook ---
$womenSheet.Columns.Item(2).ColumnWidth = 23.428571428571427
$womenSheet.Range("B8").Select()

$menSheet = $wb.Worksheets.Item("1")
$menSheet.Range("C6").Select()

$settingSheet.Select()
$settingSheet.Range("C7").Select()
